# Applies the cryptos-list refresh described by the commit diff: updated
# Price/Volume(1h) figures for most rows, plus a swap of the Uniswap/
# BitcoinCash rows (20 <-> 21).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $value, $numericLooking) {
    $cell = $ws.Range($ref)
    if ($numericLooking) {
        # A plain ".Value = <digits>" assignment would let Excel auto-convert
        # the numeric-looking string into a real Number, and pre-formatting the
        # cell as Text (or using a leading apostrophe) stamps a new NumberFormat/
        # quotePrefix style onto the cell. Neither matches the source workbook,
        # where these are plain unstyled text cells. Instead, write the text as a
        # literal-string formula, then collapse it to its static value with a
        # values-only paste so the result is a plain, unstyled text cell.
        $cell.Formula = "=""" + $value + """"
        $cell.Copy()
        $cell.PasteSpecial(-4163)
    } else {
        $cell.Value = $value
    }
}

Set-TextValue $ws "D2" "25.915.33" $false
Set-TextValue $ws "E2" "  +0.17%  " $false
Set-TextValue $ws "D3" "1.636.07" $false
Set-TextValue $ws "E3" "  +0.04%  " $false
Set-TextValue $ws "E4" "  +0.20%  " $false
Set-TextValue $ws "D5" "214.56" $true
Set-TextValue $ws "E5" "  -0.15%  " $false
Set-TextValue $ws "D6" "0.508" $true
Set-TextValue $ws "E6" "  +1.06%  " $false
Set-TextValue $ws "E7" "  +0.21%  " $false
Set-TextValue $ws "E8" "  -0.94%  " $false
Set-TextValue $ws "D10" "19.59" $true
Set-TextValue $ws "E10" "  -0.24%  " $false
Set-TextValue $ws "D11" "0.0794" $true
Set-TextValue $ws "E11" "  +0.44%  " $false
Set-TextValue $ws "D12" "1.862.48" $false
Set-TextValue $ws "D13" "4.24" $true
Set-TextValue $ws "E13" "  -0.51%  " $false
Set-TextValue $ws "D14" "1.622.58" $false
Set-TextValue $ws "E14" "  -0.92%  " $false
Set-TextValue $ws "E15" "  -1.62%  " $false
Set-TextValue $ws "E16" "  -0.48%  " $false
Set-TextValue $ws "D17" "62.55" $true
Set-TextValue $ws "E17" "  -0.55%  " $false
Set-TextValue $ws "D18" "25.931.64" $false
Set-TextValue $ws "E18" "  +0.28%  " $false
Set-TextValue $ws "E19" "  +0.16%  " $false
Set-TextValue $ws "B20" "BitcoinCash" $false
Set-TextValue $ws "C20" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch" $false
Set-TextValue $ws "D20" "193.49" $true
Set-TextValue $ws "E20" "  +1.01%  " $false
Set-TextValue $ws "B21" "Uniswap" $false
Set-TextValue $ws "C21" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni" $false
Set-TextValue $ws "D21" "4.39" $true
Set-TextValue $ws "E21" "  -1.12%  " $false
Set-TextValue $ws "E22" "  -0.52%  " $false
Set-TextValue $ws "E23" "  -0.92%  " $false
Set-TextValue $ws "E24" "  +0.24%  " $false
Set-TextValue $ws "D25" "143.73" $true
Set-TextValue $ws "E25" "  +0.91%  " $false
Set-TextValue $ws "E27" "  +2.88%  " $false
Set-TextValue $ws "E28" "  -0.14%  " $false
Set-TextValue $ws "D29" "15.43" $true
Set-TextValue $ws "E29" "  -0.55%  " $false
Set-TextValue $ws "E30" "  +0.19%  " $false
Set-TextValue $ws "E32" "  -1.32%  " $false
Set-TextValue $ws "E33" "  -0.76%  " $false
Set-TextValue $ws "E34" "  -2.66%  " $false
Set-TextValue $ws "E35" "  +1.41%  " $false
Set-TextValue $ws "D36" "0.902" $true
Set-TextValue $ws "E36" "  -0.51%  " $false
Set-TextValue $ws "D37" "1.138.96" $false
Set-TextValue $ws "E37" "  -0.82%  " $false
Set-TextValue $ws "E38" "  +0.12%  " $false
Set-TextValue $ws "E39" "  -1.17%  " $false
Set-TextValue $ws "E40" "  +0.11%  " $false
Set-TextValue $ws "E41" "  +0.14%  " $false
Set-TextValue $ws "D42" "99.47" $true
Set-TextValue $ws "E42" "  -1.02%  " $false
Set-TextValue $ws "D43" "0.797" $true
Set-TextValue $ws "E43" "  -0.63%  " $false
Set-TextValue $ws "D44" "5.43" $true
Set-TextValue $ws "E44" "  -3.53%  " $false
Set-TextValue $ws "D45" "1.771.66" $false
Set-TextValue $ws "E45" "  +0.03%  " $false
Set-TextValue $ws "E46" "  +2.95%  " $false
Set-TextValue $ws "D47" "56.36" $true
Set-TextValue $ws "E47" "  +1.33%  " $false
Set-TextValue $ws "E48" "  +3.33%  " $false
Set-TextValue $ws "E49" "  -1.20%  " $false
Set-TextValue $ws "E50" "  -0.47%  " $false
Set-TextValue $ws "D51" "7.63" $true
Set-TextValue $ws "E51" "  +0.81%  " $false

$excel.CutCopyMode = 0
